$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.986.43"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "2.301.65"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +18.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "270.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.621"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0951"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +14.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("D15").Value = "2.648.09"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.856"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "2.302.21"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").Value = "43.876.94"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.04%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  +3.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.33%  "
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0939"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0365"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.109"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.24%  "
$ws.Range("B40").Value = "MultiversX"
$ws.Range("C40").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "74.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.65%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.244"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +23.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.66%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.100"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.468"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.49%  "
